# Append a new "batch" of vaccination data (rows 87-102) to the sheet,
# mirroring a copy-down of the last existing batch (rows 71-86), and
# update the sheet view to reflect the newly scrolled/selected position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
# Column layout: A = Datenstand (only on first row of a batch), B = Datum,
# C = vacc, D = non_vacc (formula references into the previous batch for the
# overlapping days, then literal values for the brand new days), E = incidence.

$bDates = @(44502,44503,44504,44505,44506,44507,44508,44509,44510,44511,44512,44513,44514,44515)
$cVals  = @(584.20000000000005,686.1,747.2,830.5,907.4,965.6,1002.8,1076.5999999999999,1188.5999999999999,1339.4,1492.8,1604,1675,1718.3)
$eVals  = @(111.1,119.9,130.9,138.4,140.6,141.80000000000001,141.1,143.9,130,110.4,91.4,77.7,65.7,61.7)

# D formulas for rows 87-96 (overlapping days still reference the prior batch);
# rows 97-100 are brand new days with literal values.
$dFormulas = @("=D71","=D72","=D73","=D74","=D75","=D76","=D77","=D78","=D79","=D86")
$dLiterals = @(569,620.70000000000005,670.9,754.3)

for ($i = 0; $i -lt $bDates.Count; $i++) {
    $row = 87 + $i

    # B column: date, formatted like the existing date columns (copy number
    # format from the cell directly above so it reuses the same style).
    $ws.Cells.Item($row - 1, 2).Copy() | Out-Null
    $ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 2).Value = $bDates[$i]

    # C column: plain numeric value.
    $ws.Cells.Item($row, 3).Value = $cVals[$i]

    # D column: formula for the first 10 new rows, literal value afterwards.
    if ($i -lt $dFormulas.Count) {
        $ws.Cells.Item($row, 4).Formula = $dFormulas[$i]
    } else {
        $ws.Cells.Item($row, 4).Value = $dLiterals[$i - $dFormulas.Count]
    }

    # E column: plain numeric value.
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}

# Row 87 additionally carries the "Datenstand" value in column A (new batch).
$ws.Cells.Item(86, 1).Copy() | Out-Null
$ws.Cells.Item(87, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(87, 1).Value = 44515

# Trailing placeholder rows 101-102: same date-formatted style as column B,
# but left empty (no value) -- mirrors two blank rows left under the table.
$ws.Cells.Item(100, 2).Copy() | Out-Null
$ws.Cells.Item(101, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(100, 2).Copy() | Out-Null
$ws.Cells.Item(102, 2).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Sheet view: scrolled down / zoomed to show the newly appended rows --
$ws.Application.ActiveWindow.ScrollRow = 77
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.Zoom = 116
$ws.Range("E101").Select() | Out-Null
